$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1841.909
$ws.Range("I62").Value = 2045.8334
$ws.Range("J62").Value = 1597.2
$ws.Range("K62").Value = 2045.8334
$ws.Range("L62").Value = 1597.2
$ws.Range("M62").Value = -1421.8334
$ws.Range("N62").Value = -2845.2

$ws.Range("H65").Value = 1841.909
$ws.Range("I65").Value = 2045.8334
$ws.Range("J65").Value = 1597.2
$ws.Range("K65").Value = 10229.167
$ws.Range("L65").Value = 7986
$ws.Range("M65").Value = -7109.166999999999
$ws.Range("N65").Value = -14226

$ws.Range("H116").Value = 1853.125
$ws.Range("J116").Value = 2011.1111
$ws.Range("L116").Value = 2011.1111
$ws.Range("N116").Value = -8895.1111

$ws.Range("H131").Value = 739.63635
$ws.Range("I131").Value = 597.1429000000001
$ws.Range("J131").Value = 989
$ws.Range("K131").Value = 1791.4287
$ws.Range("L131").Value = 2967
$ws.Range("M131").Value = 3248.5713
$ws.Range("N131").Value = -13047

$ws.Range("H137").Value = 2780337.8
$ws.Range("I137").Value = 3573498.5
$ws.Range("J137").Value = 4275.125
$ws.Range("K137").Value = 10720495.5
$ws.Range("L137").Value = 12825.375
$ws.Range("M137").Value = -10717945.5
$ws.Range("N137").Value = -17925.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2735.7144
$ws.Range("I35").Value = 1829.6
$ws.Range("J35").Value = 5001
$ws.Range("K35").Value = 1829.6
$ws.Range("L35").Value = 5001
$ws.Range("M35").Value = -1423.6
$ws.Range("N35").Value = -5813

$ws.Range("H61").Value = 23857904
$ws.Range("I61").Value = 27054758
$ws.Range("K61").Value = 27054758
$ws.Range("M61").Value = -27054546

$ws.Range("H74").Value = 9695003
$ws.Range("I74").Value = 13212205
$ws.Range("J74").Value = 148311.42
$ws.Range("K74").Value = 13212205
$ws.Range("L74").Value = 148311.42
$ws.Range("M74").Value = -13211331
$ws.Range("N74").Value = -150059.42

$ws.Range("H77").Value = 9695003
$ws.Range("I77").Value = 13212205
$ws.Range("J77").Value = 148311.42
$ws.Range("K77").Value = 66061025
$ws.Range("L77").Value = 741557.1000000001
$ws.Range("M77").Value = -66056657
$ws.Range("N77").Value = -750293.1000000001

$ws.Range("H122").Value = 2711613.2
$ws.Range("I122").Value = 1546.8108
$ws.Range("J122").Value = 27779728
$ws.Range("K122").Value = 4640.4324
$ws.Range("L122").Value = 83339184
$ws.Range("M122").Value = -2190.4324
$ws.Range("N122").Value = -83344084

$ws.Range("H132").Value = 62912.766
$ws.Range("I132").Value = 44483.695
$ws.Range("K132").Value = 133451.085
$ws.Range("M132").Value = -130921.085

$ws.Range("H136").Value = 23857904
$ws.Range("I136").Value = 27054758
$ws.Range("K136").Value = 81164274
$ws.Range("M136").Value = -81161724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1027.6666
$ws.Range("I20").Value = 709.7778
$ws.Range("J20").Value = 1504.5
$ws.Range("K20").Value = 709.7778
$ws.Range("L20").Value = 1504.5
$ws.Range("M20").Value = -462.7778
$ws.Range("N20").Value = -1998.5

$ws.Range("H86").Value = 7796.7144
$ws.Range("I86").Value = 10773.885
$ws.Range("K86").Value = 10773.885
$ws.Range("M86").Value = -9650.885

$ws.Range("H89").Value = 7796.7144
$ws.Range("I89").Value = 10773.885
$ws.Range("K89").Value = 53869.425
$ws.Range("M89").Value = -48253.425

$ws.Range("H134").Value = 3674
$ws.Range("I134").Value = 3387.48
$ws.Range("J134").Value = 5265.778
$ws.Range("K134").Value = 10162.44
$ws.Range("L134").Value = 15797.334
$ws.Range("M134").Value = -7627.440000000001
$ws.Range("N134").Value = -20867.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 366.66666
$ws.Range("J2").Value = 400
$ws.Range("L2").Value = 400
$ws.Range("N2").Value = -626

$ws.Range("H31").Value = 3499.1428
$ws.Range("I31").Value = 2740.3076
$ws.Range("J31").Value = 3947.5454
$ws.Range("K31").Value = 2740.3076
$ws.Range("L31").Value = 3947.5454
$ws.Range("M31").Value = -2445.3076
$ws.Range("N31").Value = -4537.5454

$ws.Range("H34").Value = 3499.1428
$ws.Range("I34").Value = 2740.3076
$ws.Range("J34").Value = 3947.5454
$ws.Range("K34").Value = 2740.3076
$ws.Range("L34").Value = 3947.5454
$ws.Range("M34").Value = -2538.3076
$ws.Range("N34").Value = -4351.5454

$ws.Range("H99").Value = 4109.391
$ws.Range("I99").Value = 3480.4285
$ws.Range("J99").Value = 5087.778
$ws.Range("K99").Value = 3480.4285
$ws.Range("L99").Value = 5087.778
$ws.Range("M99").Value = -1982.4285
$ws.Range("N99").Value = -8083.778

$ws.Range("H126").Value = 4109.391
$ws.Range("I126").Value = 3480.4285
$ws.Range("J126").Value = 5087.778
$ws.Range("K126").Value = 10441.2855
$ws.Range("L126").Value = 15263.334
$ws.Range("M126").Value = -7971.2855
$ws.Range("N126").Value = -20203.334

$ws.Range("H132").Value = 33213.855
$ws.Range("I132").Value = 23496.422
$ws.Range("J132").Value = 57507.445
$ws.Range("K132").Value = 70489.266
$ws.Range("L132").Value = 172522.335
$ws.Range("M132").Value = -67959.266
$ws.Range("N132").Value = -177582.335

$ws.Range("H134").Value = 28919.834
$ws.Range("I134").Value = 2588.0334
$ws.Range("J134").Value = 94749.336
$ws.Range("K134").Value = 7764.100199999999
$ws.Range("L134").Value = 284248.008
$ws.Range("M134").Value = -5229.100199999999
$ws.Range("N134").Value = -289318.008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 611.6957
$ws.Range("I5").Value = 281.26666
$ws.Range("J5").Value = 1231.25
$ws.Range("K5").Value = 843.79998
$ws.Range("L5").Value = 3693.75
$ws.Range("M5").Value = -731.79998
$ws.Range("N5").Value = -3917.75

$ws.Range("H68").Value = 1084.0667
$ws.Range("I68").Value = 801
$ws.Range("J68").Value = 1127.6154
$ws.Range("K68").Value = 2403
$ws.Range("L68").Value = 3382.8462
$ws.Range("M68").Value = -1592
$ws.Range("N68").Value = -5004.8462

$ws.Range("H71").Value = 1084.0667
$ws.Range("I71").Value = 801
$ws.Range("J71").Value = 1127.6154
$ws.Range("K71").Value = 7209
$ws.Range("L71").Value = 10148.5386
$ws.Range("M71").Value = -3153
$ws.Range("N71").Value = -18260.5386

$ws.Range("H131").Value = 8197733
$ws.Range("I131").Value = 71428984
$ws.Range("J131").Value = 1089.1666
$ws.Range("K131").Value = 214286952
$ws.Range("L131").Value = 3267.4998
$ws.Range("M131").Value = -214281912
$ws.Range("N131").Value = -13347.4998

$ws.Range("H135").Value = 611.6957
$ws.Range("I135").Value = 281.26666
$ws.Range("J135").Value = 1231.25
$ws.Range("K135").Value = 2531.39994
$ws.Range("L135").Value = 11081.25
$ws.Range("M135").Value = 3.600059999999758
$ws.Range("N135").Value = -16151.25

$ws.Range("H137").Value = 32003.5
$ws.Range("I137").Value = 1806
$ws.Range("J137").Value = 42069.332
$ws.Range("K137").Value = 5418
$ws.Range("L137").Value = 126207.996
$ws.Range("M137").Value = -318
$ws.Range("N137").Value = -136407.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 806.6667
$ws.Range("I22").Value = 722.2222
$ws.Range("J22").Value = 933.3333
$ws.Range("K22").Value = 722.2222
$ws.Range("L22").Value = 933.3333
$ws.Range("M22").Value = -427.2222
$ws.Range("N22").Value = -1523.3333

$ws.Range("H27").Value = 806.6667
$ws.Range("I27").Value = 722.2222
$ws.Range("J27").Value = 933.3333
$ws.Range("K27").Value = 722.2222
$ws.Range("L27").Value = 933.3333
$ws.Range("M27").Value = -615.2222
$ws.Range("N27").Value = -1147.3333

$ws.Range("H68").Value = 1539
$ws.Range("I68").Value = 1492.9
$ws.Range("K68").Value = 1492.9
$ws.Range("M68").Value = -743.9000000000001

$ws.Range("H71").Value = 1539
$ws.Range("I71").Value = 1492.9
$ws.Range("K71").Value = 7464.5
$ws.Range("M71").Value = -3720.5

$ws.Range("H136").Value = 51175.188
$ws.Range("I136").Value = 31729.277
$ws.Range("K136").Value = 95187.83099999999
$ws.Range("M136").Value = -92637.83099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 135939.12
$ws.Range("I132").Value = 143501.22
$ws.Range("K132").Value = 430503.66
$ws.Range("M132").Value = -427973.66

$ws.Range("H136").Value = 46979.84
$ws.Range("I136").Value = 39312.848
$ws.Range("J136").Value = 58054.39
$ws.Range("K136").Value = 117938.544
$ws.Range("L136").Value = 174163.17
$ws.Range("M136").Value = -115388.544
$ws.Range("N136").Value = -179263.17
